# Apply the "WORKING VERSION OF LLM CLEANER THAT NEEDS TO BE SCALED UP" edit.
#
# Summary of the change (see xml diff):
#  - Row 40 gains a new "Context" cell (D40).
#  - A brand-new row is inserted at position 41 (pushing the former row 41
#    down to row 42), documenting an "Improve the loggings..." task that
#    failed.
#  - The (now) row 42 - the old row 41 - gains a new "Context" cell (D42).
#  - Two brand-new rows are appended at positions 43 and 44 describing a
#    review-quality task and a batch-processing-time measurement task.
#  - Various cosmetic window/selection settings also changed upstream; we
#    replicate what is reachable through the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlTop = -4160

# ---------------------------------------------------------------------
# 1. Insert all three new rows FIRST (before any of the surrounding cells
#    pick up column-D / column-F content), so the row-insert format
#    copy-down doesn't leave stray empty cells behind in columns that the
#    final file does not populate for these rows.
# ---------------------------------------------------------------------
$ws.Rows("41:41").Insert()
$ws.Rows("43:44").Insert()

# ---------------------------------------------------------------------
# 2. Add the missing "Context" note to existing row 40.
# ---------------------------------------------------------------------
$ws.Range("D40").Value = "Used logger.info to print the prompt per row for debugging and analysis "

# ---------------------------------------------------------------------
# 3. Fill in the brand new row 41.
# ---------------------------------------------------------------------
$ws.Range("A41").Value = "4/8/2025(Remote)"
$ws.Range("B41").Value = "Car Tracking Project"
$ws.Range("C41").Value = "Improve the loggings even more to identify where are the bottlenecks"
$ws.Range("E41").Value = "The process was going smoothly until I had to Remove all the changes that I have made today to an older working version of the branch"
$ws.Range("F41").Value = "FAILED REALLY HARD due to not testing comprhensively and waiting for the output cleaned file to see the formatting "

# E41 inherited row-40's bold "Problem" style (column E) when the row was
# inserted; the target file uses the plain left/top style there instead.
$ws.Range("E41").HorizontalAlignment = $xlLeft
$ws.Range("E41").VerticalAlignment = $xlTop
$ws.Range("E41").WrapText = $false
$ws.Range("E41").Font.Bold = $false

# ---------------------------------------------------------------------
# 4. The old row 41 ("Check if there are any errors...") is now row 42;
#    give it its new "Context" note and mark it resolved ("DONE").
# ---------------------------------------------------------------------
$ws.Range("D42").Value = "Investigated JSONDecodeError; confirmed cause was due to me writing the top_comments column name wrong by forgetting the 's' at the end of it"
$ws.Range("F42").Value = "DONE"

# ---------------------------------------------------------------------
# 5. Fill in the brand new rows 43 and 44.
# ---------------------------------------------------------------------
$ws.Range("A43").Value = "4/8/2025(Remote)"
$ws.Range("B43").Value = "Car Tracking Project"
$ws.Range("C43").Value = "Review quality of extracted content from LLM"
$ws.Range("D43").Value = "Manually reviewed several rows; evaluated clarity, correctness, and usability of output"
$ws.Range("E43").Value = "Not statisfied with the output, I feel that lots of the needed context for the diagnosis is missing from the main extracted data"
$ws.Range("F43").Value = "DONE: Prompt engineering task incoming but after finding a scalable solution for the upcoming problem"

$ws.Range("A44").Value = "4/8/2025(Remote)"
$ws.Range("B44").Value = "Car Tracking Project"
$ws.Range("C44").Value = "Measure total time taken for batch processing"
$ws.Range("F44").Value = "DONE: 44 rows took ~2 hours; extrapolated daily 700 rows to ~31.3 hours. NEEDS A SCALABLE SOLUTION INSTEAD OF `nGITHUB ACTIONS LLMS"
$ws.Range("F44").WrapText = $true
$ws.Rows("44:44").RowHeight = 28.8

# ---------------------------------------------------------------------
# 6. Cosmetic view-state updates (best effort through the object model).
# ---------------------------------------------------------------------
$ws.Range("D48").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 4

$excel.ActiveWindow.Left = 11424
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 11712
$excel.ActiveWindow.Height = 12336
